# Daily attendance processing - 2025-12-26 22:57:32
#
# For every data row in the "Session Analysis Results" sheet, the
# "Recorded By" column (G) lists the user(s) who recorded that session,
# separated by ", ". Rows that list exactly two names had the order of
# those two names swapped (this does not apply to the "backup@backdoor.com"
# audit rows, whose ordering is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$changed = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ($val -eq "Recorded By") {
        continue
    }

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    if ($val -match "backup@backdoor.com") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -eq 2) {
        $swapped = $parts[1] + ", " + $parts[0]
        $cell.Value = $swapped
        $changed = $changed + 1
    }
}

Write-Output "Swapped Recorded By order on $changed rows"
